# Insert a new weekly price record at row 433 (Vega Modelo de Temuco - Perejil),
# shifting the existing rows 433:466 down to 434:467.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(433).Insert()

$ws.Range("A433").Value = 10
$ws.Range("B433").Value = "Vega Modelo de Temuco"
$ws.Range("C433").Value = "La Araucanía"
$ws.Range("D433").Value = 45021
$ws.Range("E433").Value = 9
$ws.Range("F433").Value = 100112044
$ws.Range("G433").Value = "Perejil"
$ws.Range("H433").Value = "Sin especificar"
$ws.Range("I433").Value = "Primera"
$ws.Range("J433").Value = 55
$ws.Range("K433").Value = 4000
$ws.Range("L433").Value = 4000
$ws.Range("M433").Value = 4000
$ws.Range("N433").Value = "$/docena de atados (3 kilos)"
$ws.Range("O433").Value = "Provincia de Cautín"
$ws.Range("P433").Value = 1333
$ws.Range("Q433").Value = 3
$ws.Range("R433").Value = "Hortaliza"
